# This workbook's data rows (2..40) represent daily price records. The
# commit "Fruta / hortaliza, semanal" re-shuffles which date/price/origin
# block sits on which row (the rows 2..40 are a permutation of each
# other for columns D, I, J, K, L, M, O, P - all other columns are
# unchanged). We snapshot the current ("before") values for those
# columns, then write them back out according to the row permutation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column numbers: D=4, I=9, J=10, K=11, L=12, M=13, O=15, P=16
$firstRow = 2
$lastRow = 40

# Snapshot all current values first, since rows will be overwritten.
$dVals = @{}
$iVals = @{}
$jVals = @{}
$kVals = @{}
$lVals = @{}
$mVals = @{}
$oVals = @{}
$pVals = @{}

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $dVals[$r] = $ws.Cells.Item($r, 4).Value2
    $iVals[$r] = $ws.Cells.Item($r, 9).Value2
    $jVals[$r] = $ws.Cells.Item($r, 10).Value2
    $kVals[$r] = $ws.Cells.Item($r, 11).Value2
    $lVals[$r] = $ws.Cells.Item($r, 12).Value2
    $mVals[$r] = $ws.Cells.Item($r, 13).Value2
    $oVals[$r] = $ws.Cells.Item($r, 15).Value2
    $pVals[$r] = $ws.Cells.Item($r, 16).Value2
}

# Destination row -> source row (data formerly on the source row now
# appears on the destination row).
$rowMap = @{
    2 = 5;   3 = 10;  4 = 19;  5 = 36;  6 = 8;   7 = 2;   8 = 12;  9 = 14;
    10 = 35; 11 = 37; 12 = 27; 13 = 39; 14 = 22; 15 = 20; 16 = 26; 17 = 40;
    18 = 38; 19 = 6;  20 = 29; 21 = 15; 22 = 4;  23 = 13; 24 = 21; 25 = 18;
    26 = 33; 27 = 9;  28 = 25; 29 = 23; 30 = 24; 31 = 28; 32 = 32; 33 = 30;
    34 = 31; 35 = 7;  36 = 16; 37 = 17; 38 = 11; 39 = 34; 40 = 3
}

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $src = $rowMap[$r]

    $ws.Cells.Item($r, 4).Value = $dVals[$src]
    $ws.Cells.Item($r, 9).Value = $iVals[$src]
    $ws.Cells.Item($r, 10).Value = $jVals[$src]
    $ws.Cells.Item($r, 11).Value = $kVals[$src]
    $ws.Cells.Item($r, 12).Value = $lVals[$src]
    $ws.Cells.Item($r, 13).Value = $mVals[$src]
    $ws.Cells.Item($r, 15).Value = $oVals[$src]
    $ws.Cells.Item($r, 16).Value = $pVals[$src]
}
